$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Sectm1a"
$ws.Range("C2").Value = "Cd7"
$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1005446666666667
$ws.Range("H2").Value = 0.301634
$ws.Range("I2").Value = 0.5436002609563836
$ws.Range("J2").Value = 0.5436002609563836
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3679763333333333
$ws.Range("N2").Value = 1.103929
$ws.Range("O2").Value = 0.4187915190780531
$ws.Range("P2").Value = 0.4187915190780531
$ws.Range("Q2").Value = 0.03699805777622222
$ws.Range("R2").Value = 0.332982519986
$ws.Range("S2").Value = 0.2276551790571499
$ws.Range("T2").Value = 0.2276551790571499

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Sectm1a"
$ws.Range("C3").Value = "Cd7"
$ws.Range("D3").Value = "Neutrophils"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1005446666666667
$ws.Range("H3").Value = 0.301634
$ws.Range("I3").Value = 0.5436002609563836
$ws.Range("J3").Value = 0.5436002609563836
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.371976
$ws.Range("N3").Value = 1.115928
$ws.Range("O3").Value = 0.4233435142130823
$ws.Range("P3").Value = 0.4233435142130823
$ws.Range("Q3").Value = 0.037400202928
$ws.Range("R3").Value = 0.336601826352
$ws.Range("S3").Value = 0.230129644800424
$ws.Range("T3").Value = 0.230129644800424

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Sectm1a"
$ws.Range("C4").Value = "Cd7"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1005446666666667
$ws.Range("H4").Value = 0.301634
$ws.Range("I4").Value = 0.5436002609563836
$ws.Range("J4").Value = 0.5436002609563836
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.13871
$ws.Range("N4").Value = 0.41613
$ws.Range("O4").Value = 0.1578649667088647
$ws.Range("P4").Value = 0.1578649667088647
$ws.Range("Q4").Value = 0.01394655071333333
$ws.Range("R4").Value = 0.12551895642
$ws.Range("S4").Value = 0.08581543709880964
$ws.Range("T4").Value = 0.08581543709880964

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Sectm1a"
$ws.Range("C5").Value = "Cd7"
$ws.Range("D5").Value = "Inflammatory-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.03475566666666666
$ws.Range("H5").Value = 0.104267
$ws.Range("I5").Value = 0.1879084201686124
$ws.Range("J5").Value = 0.1879084201686124
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.3679763333333333
$ws.Range("N5").Value = 1.103929
$ws.Range("O5").Value = 0.4187915190780531
$ws.Range("P5").Value = 0.4187915190780531
$ws.Range("Q5").Value = 0.01278926278255555
$ws.Range("R5").Value = 0.115103365043
$ws.Range("S5").Value = 0.07869445272997026
$ws.Range("T5").Value = 0.07869445272997026

# Row 6
$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Sectm1a"
$ws.Range("C6").Value = "Cd7"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.03475566666666666
$ws.Range("H6").Value = 0.104267
$ws.Range("I6").Value = 0.1879084201686124
$ws.Range("J6").Value = 0.1879084201686124
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.371976
$ws.Range("N6").Value = 1.115928
$ws.Range("O6").Value = 0.4233435142130823
$ws.Range("P6").Value = 0.4233435142130823
$ws.Range("Q6").Value = 0.012928273864
$ws.Range("R6").Value = 0.116354464776
$ws.Range("S6").Value = 0.0795498109444088
$ws.Range("T6").Value = 0.0795498109444088

# Row 7
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Sectm1a"
$ws.Range("C7").Value = "Cd7"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.03475566666666666
$ws.Range("H7").Value = 0.104267
$ws.Range("I7").Value = 0.1879084201686124
$ws.Range("J7").Value = 0.1879084201686124
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.13871
$ws.Range("N7").Value = 0.41613
$ws.Range("O7").Value = 0.1578649667088647
$ws.Range("P7").Value = 0.1578649667088647
$ws.Range("Q7").Value = 0.004820958523333333
$ws.Range("R7").Value = 0.04338862671
$ws.Range("S7").Value = 0.02966415649423335
$ws.Range("T7").Value = 0.02966415649423335

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Sectm1a"
$ws.Range("C8").Value = "Cd7"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.04966033333333333
$ws.Range("H8").Value = 0.148981
$ws.Range("I8").Value = 0.268491318875004
$ws.Range("J8").Value = 0.268491318875004
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.3679763333333333
$ws.Range("N8").Value = 1.103929
$ws.Range("O8").Value = 0.4187915190780531
$ws.Range("P8").Value = 0.4187915190780531
$ws.Range("Q8").Value = 0.01827382737211111
$ws.Range("R8").Value = 0.164464446349
$ws.Range("S8").Value = 0.1124418872909329
$ws.Range("T8").Value = 0.1124418872909329

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Sectm1a"
$ws.Range("C9").Value = "Cd7"
$ws.Range("D9").Value = "Neutrophils"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.04966033333333333
$ws.Range("H9").Value = 0.148981
$ws.Range("I9").Value = 0.268491318875004
$ws.Range("J9").Value = 0.268491318875004
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.371976
$ws.Range("N9").Value = 1.115928
$ws.Range("O9").Value = 0.4233435142130823
$ws.Range("P9").Value = 0.4233435142130823
$ws.Range("Q9").Value = 0.018472452152
$ws.Range("R9").Value = 0.166252069368
$ws.Range("S9").Value = 0.1136640584682495
$ws.Range("T9").Value = 0.1136640584682495

# Row 10
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Sectm1a"
$ws.Range("C10").Value = "Cd7"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.04966033333333333
$ws.Range("H10").Value = 0.148981
$ws.Range("I10").Value = 0.268491318875004
$ws.Range("J10").Value = 0.268491318875004
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.13871
$ws.Range("N10").Value = 0.41613
$ws.Range("O10").Value = 0.1578649667088647
$ws.Range("P10").Value = 0.1578649667088647
$ws.Range("Q10").Value = 0.006888384836666667
$ws.Range("R10").Value = 0.06199546353
$ws.Range("S10").Value = 0.04238537311582168
$ws.Range("T10").Value = 0.04238537311582168
